$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 stay text-formatted so numeric-looking strings
# (e.g. "0.650", "1.00", "69.261.31") are not coerced to numbers
# and keep their exact original formatting, matching the source file
# where these cells are stored as inline strings (t="inlineStr").
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '69.261.31'
$ws.Range("E2").Value = '  -0.23%  '

# Row 3
$ws.Range("D3").Value = '3.679.77'
$ws.Range("E3").Value = '  -0.11%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '681.55'
$ws.Range("E5").Value = '  -0.20%  '

# Row 6
$ws.Range("D6").Value = '158.32'
$ws.Range("E6").Value = '  -2.41%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '0.494'
$ws.Range("E8").Value = '  -0.80%  '

# Row 9
$ws.Range("E9").Value = '  -0.99%  '

# Row 10
$ws.Range("D10").Value = '7.04'
$ws.Range("E10").Value = '  -2.97%  '

# Row 11
$ws.Range("E11").Value = '  -1.69%  '

# Row 12
$ws.Range("D12").Value = '0.0000233'
$ws.Range("E12").Value = '  -1.92%  '

# Row 13
$ws.Range("D13").Value = '4.301.43'
$ws.Range("E13").Value = '  +0.04%  '

# Row 14
$ws.Range("D14").Value = '32.30'
$ws.Range("E14").Value = '  -3.31%  '

# Row 15
$ws.Range("D15").Value = '3.672.27'
$ws.Range("E15").Value = '  -0.10%  '

# Row 16
$ws.Range("D16").Value = '69.261.74'
$ws.Range("E16").Value = '  -0.27%  '

# Row 17
$ws.Range("E17").Value = '  +1.98%  '

# Row 18
$ws.Range("E18").Value = '  -1.67%  '

# Row 19
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  -3.43%  '

# Row 20
$ws.Range("D20").Value = '472.36'
$ws.Range("E20").Value = '  -0.08%  '

# Row 21
$ws.Range("D21").Value = '10.06'
$ws.Range("E21").Value = '  +2.33%  '

# Row 22
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").Value = '  -2.21%  '

# Row 23
$ws.Range("D23").Value = '79.92'
$ws.Range("E23").Value = '  -0.12%  '

# Row 24
$ws.Range("D24").Value = '3.826.00'

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").Value = '0.0000122'
$ws.Range("E26").Value = '  -5.20%  '

# Row 27
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -4.33%  '

# Row 28
$ws.Range("D28").Value = '9.10'
$ws.Range("E28").Value = '  -4.95%  '

# Row 29
$ws.Range("D29").Value = '2.70'
$ws.Range("E29").Value = '  -1.46%  '

# Row 30
$ws.Range("E30").Value = '  -5.01%  '

# Row 31
$ws.Range("D31").Value = '6.63'
$ws.Range("E31").Value = '  -2.84%  '

# Row 32
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '26.87'
$ws.Range("E33").Value = '  -0.35%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '1.98'
$ws.Range("E34").Value = '  -5.31%  '

# Row 35
$ws.Range("D35").Value = '3.659.34'
$ws.Range("E35").Value = '  +0.45%  '

# Row 36
$ws.Range("D36").Value = '0.158'
$ws.Range("E36").Value = '  -4.73%  '

# Row 37
$ws.Range("E37").Value = '  -3.47%  '

# Row 38
$ws.Range("E38").Value = '  -1.10%  '

# Row 40
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  +4.03%  '

# Row 41
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '0.0904'
$ws.Range("E42").Value = '  -3.48%  '

# Row 43
$ws.Range("D43").Value = '170.97'
$ws.Range("E43").Value = '  +9.54%  '

# Row 44
$ws.Range("E44").Value = '  -1.50%  '

# Row 45
$ws.Range("D45").Value = '47.49'
$ws.Range("E45").Value = '  -1.65%  '

# Row 46
$ws.Range("D46").Value = '2.72'
$ws.Range("E46").Value = '  -6.07%  '

# Row 47
$ws.Range("D47").Value = '0.000280'
$ws.Range("E47").Value = '  -3.48%  '

# Row 48
$ws.Range("D48").Value = '1.10'
$ws.Range("E48").Value = '  +2.14%  '

# Row 49
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  -5.40%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '26.98'
$ws.Range("E50").Value = '  -3.11%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '7.76'
$ws.Range("E51").Value = '  -3.66%  '

